# Updated sceneration output to simulation input method
# Rewrites rows 2-11 (simulation results: arrival/departure time & SoC,
# vehicle model, battery capacity, charging power) with new sample data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 44348.41666666666
$ws.Cells.Item(2, 3).Value = 44349.46875
$ws.Cells.Item(2, 4).Value = 0.44
$ws.Cells.Item(2, 5).Value = 0.861
$ws.Cells.Item(2, 6).Value = "kia niro ev"
$ws.Cells.Item(2, 7).Value = 64.8
$ws.Cells.Item(2, 8).Value = 11
$ws.Cells.Item(2, 9).Value = 80

# Row 3
$ws.Cells.Item(3, 2).Value = 44348.41666666666
$ws.Cells.Item(3, 3).Value = 44348.69791666666
$ws.Cells.Item(3, 4).Value = 0.425
$ws.Cells.Item(3, 5).Value = 0.658
$ws.Cells.Item(3, 6).Value = "kia niro ev"
$ws.Cells.Item(3, 7).Value = 64.8
$ws.Cells.Item(3, 8).Value = 11
$ws.Cells.Item(3, 9).Value = 80

# Row 4
$ws.Cells.Item(4, 2).Value = 44348.35416666666
$ws.Cells.Item(4, 3).Value = 44349.6875
$ws.Cells.Item(4, 4).Value = 0.401
$ws.Cells.Item(4, 5).Value = 0.802
$ws.Cells.Item(4, 6).Value = "audi e-tron gt rs"
$ws.Cells.Item(4, 7).Value = 85
$ws.Cells.Item(4, 8).Value = 11
$ws.Cells.Item(4, 9).Value = 265

# Row 5
$ws.Cells.Item(5, 2).Value = 44348.4375
$ws.Cells.Item(5, 3).Value = 44349.71875
$ws.Cells.Item(5, 4).Value = 0.495
$ws.Cells.Item(5, 5).Value = 0.8470000000000001
$ws.Cells.Item(5, 6).Value = "tesla y long range"
$ws.Cells.Item(5, 7).Value = 75
$ws.Cells.Item(5, 8).Value = 11
$ws.Cells.Item(5, 9).Value = 250

# Row 6
$ws.Cells.Item(6, 2).Value = 44348.58333333334
$ws.Cells.Item(6, 3).Value = 44349.76041666666
$ws.Cells.Item(6, 4).Value = 0.377
$ws.Cells.Item(6, 5).Value = 0.672
$ws.Cells.Item(6, 6).Value = "renault megane e-tech ev60 220hp"
$ws.Cells.Item(6, 7).Value = 60
$ws.Cells.Item(6, 8).Value = 22
$ws.Cells.Item(6, 9).Value = 130

# Row 7
$ws.Cells.Item(7, 2).Value = 44349.48958333334
$ws.Cells.Item(7, 3).Value = 44350
$ws.Cells.Item(7, 4).Value = 0.421
$ws.Cells.Item(7, 5).Value = 0.628
$ws.Cells.Item(7, 6).Value = "kia ev6 gt"
$ws.Cells.Item(7, 7).Value = 74
$ws.Cells.Item(7, 8).Value = 11
$ws.Cells.Item(7, 9).Value = 233

# Row 8
$ws.Cells.Item(8, 2).Value = 44349.34375
$ws.Cells.Item(8, 3).Value = 44349.77083333334
$ws.Cells.Item(8, 4).Value = 0.5940000000000001
$ws.Cells.Item(8, 5).Value = 0.8390000000000001
$ws.Cells.Item(8, 6).Value = "bmw i4 edrive40"
$ws.Cells.Item(8, 7).Value = 80.7
$ws.Cells.Item(8, 8).Value = 11
$ws.Cells.Item(8, 9).Value = 200

# Row 9
$ws.Cells.Item(9, 2).Value = 44349.35416666666
$ws.Cells.Item(9, 3).Value = 44349.86458333334
$ws.Cells.Item(9, 4).Value = 0.36
$ws.Cells.Item(9, 5).Value = 0.807
$ws.Cells.Item(9, 6).Value = "mercedes eqs 580 4matic"
$ws.Cells.Item(9, 7).Value = 107.8
$ws.Cells.Item(9, 8).Value = 11
$ws.Cells.Item(9, 9).Value = 207

# Row 10
$ws.Cells.Item(10, 2).Value = 44349.36458333334
$ws.Cells.Item(10, 3).Value = 44350
$ws.Cells.Item(10, 4).Value = 0.492
$ws.Cells.Item(10, 5).Value = 0.8320000000000001
$ws.Cells.Item(10, 6).Value = "mercedes eqe 350+"
$ws.Cells.Item(10, 7).Value = 90.6
$ws.Cells.Item(10, 8).Value = 11
$ws.Cells.Item(10, 9).Value = 170

# Row 11
$ws.Cells.Item(11, 2).Value = 44349.03125
$ws.Cells.Item(11, 3).Value = 44350
$ws.Cells.Item(11, 4).Value = 0.438
$ws.Cells.Item(11, 5).Value = 0.5860000000000001
$ws.Cells.Item(11, 6).Value = "bmw i4 m50"
$ws.Cells.Item(11, 7).Value = 80.7
$ws.Cells.Item(11, 8).Value = 11
$ws.Cells.Item(11, 9).Value = 200

